$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.965.77"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.642.63"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.872.67"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.631.59"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "27.932.65"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "1.473.35"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.926"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.24%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "1.782.27"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
